$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.438681721687317
$ws.Range("B1").Value = 3.897445678710938
$ws.Range("C1").Value = 3.072275876998901
$ws.Range("D1").Value = 2.469514131546021
$ws.Range("E1").Value = 1.400667309761047
